$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $origStyle = $cell.Style
    $cell.Value = "'" + $text
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "61.420.98"
Set-TextValue $ws.Range("E2") "  +1.01%  "
Set-TextValue $ws.Range("D3") "2.378.98"
Set-TextValue $ws.Range("E3") "  +0.96%  "
Set-TextValue $ws.Range("E4") "  +0.00%  "
Set-TextValue $ws.Range("D5") "553.04"
Set-TextValue $ws.Range("E5") "  +2.56%  "
Set-TextValue $ws.Range("D6") "139.73"
Set-TextValue $ws.Range("E6") "  +1.61%  "
Set-TextValue $ws.Range("E7") "  +0.01%  "
Set-TextValue $ws.Range("D8") "0.525"
Set-TextValue $ws.Range("E8") "  +0.78%  "
Set-TextValue $ws.Range("D9") "2.380.70"
Set-TextValue $ws.Range("E9") "  +1.10%  "
Set-TextValue $ws.Range("E10") "  +3.37%  "
Set-TextValue $ws.Range("E11") "  +2.28%  "
Set-TextValue $ws.Range("D12") "5.37"
Set-TextValue $ws.Range("E12") "  +2.37%  "
Set-TextValue $ws.Range("E13") "  +3.28%  "
Set-TextValue $ws.Range("D14") "25.56"
Set-TextValue $ws.Range("E14") "  +3.11%  "
Set-TextValue $ws.Range("D15") "0.0000173"
Set-TextValue $ws.Range("E15") "  +6.69%  "
Set-TextValue $ws.Range("D16") "2.806.98"
Set-TextValue $ws.Range("E16") "  +0.98%  "
Set-TextValue $ws.Range("D17") "61.249.60"
Set-TextValue $ws.Range("E17") "  +0.99%  "
Set-TextValue $ws.Range("D18") "2.377.72"
Set-TextValue $ws.Range("E18") "  +1.02%  "
Set-TextValue $ws.Range("D19") "10.95"
Set-TextValue $ws.Range("E19") "  +3.53%  "
Set-TextValue $ws.Range("D20") "4.16"
Set-TextValue $ws.Range("E20") "  +2.67%  "
Set-TextValue $ws.Range("D21") "320.72"
Set-TextValue $ws.Range("E21") "  +1.75%  "
Set-TextValue $ws.Range("E22") "  +1.70%  "
Set-TextValue $ws.Range("E23") "  -0.08%  "
Set-TextValue $ws.Range("D24") "64.30"
Set-TextValue $ws.Range("E24") "  +1.77%  "
Set-TextValue $ws.Range("E25") "  -8.02%  "
Set-TextValue $ws.Range("D26") "8.88"
Set-TextValue $ws.Range("E26") "  +4.97%  "
Set-TextValue $ws.Range("E27") "  -0.15%  "
Set-TextValue $ws.Range("D28") "2.494.87"
Set-TextValue $ws.Range("E28") "  +0.98%  "
Set-TextValue $ws.Range("D29") "8.18"
Set-TextValue $ws.Range("E29") "  +2.83%  "
Set-TextValue $ws.Range("D30") "521.52"
Set-TextValue $ws.Range("E30") "  +3.60%  "
Set-TextValue $ws.Range("D31") "0.0₃0905"
Set-TextValue $ws.Range("E31") "  +1.21%  "
Set-TextValue $ws.Range("D32") "1.40"
Set-TextValue $ws.Range("E32") "  +0.87%  "
Set-TextValue $ws.Range("D33") "0.148"
Set-TextValue $ws.Range("E33") "  +2.31%  "
Set-TextValue $ws.Range("E34") "  +3.18%  "
Set-TextValue $ws.Range("D35") "1.51"
Set-TextValue $ws.Range("E35") "  -0.87%  "
Set-TextValue $ws.Range("E36") "  +0.01%  "
Set-TextValue $ws.Range("D37") "5.54"
Set-TextValue $ws.Range("E37") "  +5.58%  "
Set-TextValue $ws.Range("D38") "4.70"
Set-TextValue $ws.Range("E38") "  +2.96%  "
Set-TextValue $ws.Range("E39") "  +6.12%  "
Set-TextValue $ws.Range("E40") "  +1.74%  "
Set-TextValue $ws.Range("D41") "18.53"
Set-TextValue $ws.Range("E41") "  +0.11%  "
Set-TextValue $ws.Range("D42") "146.55"
Set-TextValue $ws.Range("E42") "  +5.71%  "
Set-TextValue $ws.Range("E43") "  -0.04%  "
Set-TextValue $ws.Range("D44") "41.36"
Set-TextValue $ws.Range("E44") "  +3.10%  "
Set-TextValue $ws.Range("D45") "147.40"
Set-TextValue $ws.Range("E45") "  +6.34%  "
Set-TextValue $ws.Range("E46") "  +2.32%  "
Set-TextValue $ws.Range("E47") "  +2.71%  "
Set-TextValue $ws.Range("D48") "0.0523"
Set-TextValue $ws.Range("E48") "  +2.48%  "
Set-TextValue $ws.Range("D49") "19.81"
Set-TextValue $ws.Range("E49") "  +1.60%  "
Set-TextValue $ws.Range("D50") "0.582"
Set-TextValue $ws.Range("E50") "  +2.52%  "
Set-TextValue $ws.Range("E51") "  +1.42%  "
